$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs","Slit2","Gpc1","ECs",3,1,0.1645843333333333,0.493753,0.03485847193389392,0.03485847193389392,2,0.6666666666666666,1.628421,4.885263,0.048329411442081,0.048329411442081,0.268012584671,2.412113262039,0.001684689432335392,0.001684689432335392),
    @("ECs","Slit2","Gpc1","FAPs",3,1,0.1645843333333333,0.493753,0.03485847193389392,0.03485847193389392,3,1,9.459065000000001,28.377195,0.2807327123897247,0.2807327123897247,1.556813906981667,14.011325162835,0.009785913375763133,0.009785913375763131),
    @("ECs","Slit2","Gpc1","M2",3,1,0.1645843333333333,0.493753,0.03485847193389392,0.03485847193389392,3,1,0.6418243333333334,1.925473,0.01904850912583786,0.01904850912583786,0.1056342300187778,0.9507080701690001,0.0006640019207455411,0.0006640019207455411),
    @("ECs","Slit2","Gpc1","sCs",3,1,0.1645843333333333,0.493753,0.03485847193389392,0.03485847193389392,3,1,21.96489266666667,65.894678,0.6518893670423564,0.6518893670423563,3.615077216281556,32.535694946534,0.02272386720504985,0.02272386720504985),
    @("FAPs","Slit2","Gpc1","ECs",3,1,3.368329,10.104987,0.7134020567608963,0.7134020567608964,2,0.6666666666666666,1.628421,4.885263,0.048329411442081,0.048329411442081,5.485057678509,49.365519106581,0.03447830152482418,0.03447830152482419),
    @("FAPs","Slit2","Gpc1","FAPs",3,1,3.368329,10.104987,0.7134020567608963,0.7134020567608964,3,1,9.459065000000001,28.377195,0.2807327123897247,0.2807327123897247,31.861242952385,286.751186571465,0.2002752944188948,0.2002752944188948),
    @("FAPs","Slit2","Gpc1","M2",3,1,3.368329,10.104987,0.7134020567608963,0.7134020567608964,3,1,0.6418243333333334,1.925473,0.01904850912583786,0.01904850912583786,2.161875514872333,19.456879633851,0.01358924558860143,0.01358924558860143),
    @("FAPs","Slit2","Gpc1","sCs",3,1,3.368329,10.104987,0.7134020567608963,0.7134020567608964,3,1,21.96489266666667,65.894678,0.6518893670423564,0.6518893670423563,73.98498495102066,665.864864559186,0.465059215228576,0.4650592152285759),
    @("sCs","Slit2","Gpc1","ECs",3,1,1.188588333333333,3.565765,0.2517394713052097,0.2517394713052098,2,0.6666666666666666,1.628421,4.885263,0.048329411442081,0.048329411442081,1.935522202355,17.419699821195,0.01216642048492142,0.01216642048492143),
    @("sCs","Slit2","Gpc1","FAPs",3,1,1.188588333333333,3.565765,0.2517394713052097,0.2517394713052098,3,1,9.459065000000001,28.377195,0.2807327123897247,0.2807327123897247,11.24293430324167,101.186408729175,0.07067150459506681,0.07067150459506681),
    @("sCs","Slit2","Gpc1","M2",3,1,1.188588333333333,3.565765,0.2517394713052097,0.2517394713052098,3,1,0.6418243333333334,1.925473,0.01904850912583786,0.01904850912583786,0.7628649146494446,6.865784231845001,0.004795261616490886,0.004795261616490886),
    @("sCs","Slit2","Gpc1","sCs",3,1,1.188588333333333,3.565765,0.2517394713052097,0.2517394713052098,3,1,21.96489266666667,65.894678,0.6518893670423564,0.6518893670423563,26.10721516651889,234.96493649867,0.1641062846087306,0.1641062846087306),
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
